$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.52286145407815
$ws.Range("C2").Value = 8.125618589812118
$ws.Range("E2").Value = 14.38647635573963
$ws.Range("F2").Value = 43.34187850778581
$ws.Range("G2").Value = 46.43838137436528
$ws.Range("H2").Value = 18.74069898220959
$ws.Range("J2").Value = 9.10876621359731
$ws.Range("K2").Value = 12.22153038638588
$ws.Range("L2").Value = 11.88983048196474
$ws.Range("N2").Value = 21.68470469119048
$ws.Range("B3").Value = 16.34650949619729
$ws.Range("C3").Value = 8.099939141803356
$ws.Range("E3").Value = 14.38743306750946
$ws.Range("F3").Value = 43.33592584803153
$ws.Range("G3").Value = 46.47913914430563
$ws.Range("H3").Value = 18.7866567260709
$ws.Range("J3").Value = 9.117806293246066
$ws.Range("K3").Value = 12.10248626522043
$ws.Range("L3").Value = 11.87879153295757
$ws.Range("N3").Value = 21.74958963312509
$ws.Range("B4").Value = 16.24104519152755
$ws.Range("C4").Value = 8.083796214296605
$ws.Range("E4").Value = 14.39002244022426
$ws.Range("F4").Value = 43.34155791407336
$ws.Range("G4").Value = 46.515602263716
$ws.Range("H4").Value = 18.81775001090841
$ws.Range("J4").Value = 9.123627386836754
$ws.Range("K4").Value = 12.0314687895017
$ws.Range("L4").Value = 11.8737695242833
$ws.Range("N4").Value = 21.7913230628302
$ws.Range("B5").Value = 16.19882342886564
$ws.Range("C5").Value = 8.077123981420712
$ws.Range("E5").Value = 14.3915824004535
$ws.Range("F5").Value = 43.34619002035595
$ws.Range("G5").Value = 46.53333145196748
$ws.Range("H5").Value = 18.83114323127051
$ws.Range("J5").Value = 9.126067757619735
$ws.Range("K5").Value = 12.0030809905039
$ws.Range("L5").Value = 11.87216674421864
$ws.Range("N5").Value = 21.80880724686681
$ws.Range("B6").Value = 16.19185954649728
$ws.Range("C6").Value = 8.076010410826337
$ws.Range("E6").Value = 14.39187195433515
$ws.Range("F6").Value = 43.34710035129699
$ws.Range("G6").Value = 46.53644852246196
$ws.Range("H6").Value = 18.83341078461773
$ws.Range("J6").Value = 9.126477106631341
$ws.Range("K6").Value = 11.99840144761394
$ws.Range("L6").Value = 11.87192746337482
$ws.Range("N6").Value = 21.81173935752937
$ws.Range("B7").Value = 16.24047265214323
$ws.Range("C7").Value = 8.08370660946491
$ws.Range("E7").Value = 14.39004143290559
$ws.Range("F7").Value = 43.34161092062831
$ws.Range("G7").Value = 46.51582975341589
$ws.Range("H7").Value = 18.81792771207673
$ws.Range("J7").Value = 9.123660021988615
$ws.Range("K7").Value = 12.03108366577832
$ws.Range("L7").Value = 11.87374610926718
$ws.Range("N7").Value = 21.79155692572964
$ws.Range("B8").Value = 16.46149969747889
$ws.Range("C8").Value = 8.11684263734967
$ws.Range("E8").Value = 14.38639158652842
$ws.Range("F8").Value = 43.33789989349095
$ws.Range("G8").Value = 46.45005757990155
$ws.Range("H8").Value = 18.75594820559559
$ws.Range("J8").Value = 9.111827235413173
$ws.Range("K8").Value = 12.18007159727758
$ws.Range("L8").Value = 11.88566117569304
$ws.Range("N8").Value = 21.70668471885461
$ws.Range("B9").Value = 16.91512832613881
$ws.Range("C9").Value = 8.178817020822278
$ws.Range("E9").Value = 14.39505405043738
$ws.Range("F9").Value = 43.40418677267964
$ws.Range("G9").Value = 46.41204291836164
$ws.Range("H9").Value = 18.65724058365522
$ws.Range("J9").Value = 9.090758846944178
$ws.Range("K9").Value = 12.48733019775341
$ws.Range("L9").Value = 11.92285546666409
$ws.Range("N9").Value = 21.55521864788574
$ws.Range("B10").Value = 17.25777452549211
$ws.Range("C10").Value = 8.222486798300455
$ws.Range("E10").Value = 14.41096844715387
$ws.Range("F10").Value = 43.49749442718206
$ws.Range("G10").Value = 46.4397774896517
$ws.Range("H10").Value = 18.59866936093044
$ws.Range("J10").Value = 9.076567411853524
$ws.Range("K10").Value = 12.72037743147844
$ws.Range("L10").Value = 11.95846175505744
$ws.Range("N10").Value = 21.45297805889639
$ws.Range("B11").Value = 17.4150517600306
$ws.Range("C11").Value = 8.241942335715571
$ws.Range("E11").Value = 14.42025871463958
$ws.Range("F11").Value = 43.54954965917661
$ws.Range("G11").Value = 46.46449656649488
$ws.Range("H11").Value = 18.57505838884792
$ws.Range("J11").Value = 9.070387840549353
$ws.Range("K11").Value = 12.82757057142612
$ws.Range("L11").Value = 11.97642116387141
$ws.Range("N11").Value = 21.40841170562017
$ws.Range("B12").Value = 17.47475764187111
$ws.Range("C12").Value = 8.249250050219793
$ws.Range("E12").Value = 14.42406929908644
$ws.Range("F12").Value = 43.57063433704084
$ws.Range("G12").Value = 46.47559587007488
$ws.Range("H12").Value = 18.56655413684907
$ws.Range("J12").Value = 9.068087280073602
$ws.Range("K12").Value = 12.86829662191687
$ws.Range("L12").Value = 11.98347178817919
$ws.Range("N12").Value = 21.3918137541656
$ws.Range("B13").Value = 17.46189316107349
$ws.Range("C13").Value = 8.247678868533898
$ws.Range("E13").Value = 14.42323564871137
$ws.Range("F13").Value = 43.56603249504315
$ws.Range("G13").Value = 46.47312814828883
$ws.Range("H13").Value = 18.56836624909663
$ws.Range("J13").Value = 9.068580993021936
$ws.Range("K13").Value = 12.85952012158659
$ws.Range("L13").Value = 11.98194226007511
$ws.Range("N13").Value = 21.39537605723535
$ws.Range("B14").Value = 17.419961130383
$ws.Range("C14").Value = 8.242544745834387
$ws.Range("E14").Value = 14.4205663642972
$ws.Range("F14").Value = 43.55125685964919
$ws.Range("G14").Value = 46.46537488095957
$ws.Range("H14").Value = 18.57434998378625
$ws.Range("J14").Value = 9.070197781408826
$ws.Range("K14").Value = 12.83091863661701
$ws.Range("L14").Value = 11.9769962381622
$ws.Range("N14").Value = 21.40704060967877
$ws.Range("B15").Value = 17.39429429222531
$ws.Range("C15").Value = 8.239392154882909
$ws.Range("E15").Value = 14.41896937861712
$ws.Range("F15").Value = 43.54238479600085
$ws.Range("G15").Value = 46.46085214006357
$ws.Range("H15").Value = 18.57807208124978
$ws.Range("J15").Value = 9.071193250137586
$ws.Range("K15").Value = 12.81341585486541
$ws.Range("L15").Value = 11.97399907652178
$ws.Range("N15").Value = 21.41422170306085
$ws.Range("B16").Value = 17.24751977032933
$ws.Range("C16").Value = 8.221207021130885
$ws.Range("E16").Value = 14.41040237692073
$ws.Range("F16").Value = 43.49428509292749
$ws.Range("G16").Value = 46.43840547206575
$ws.Range("H16").Value = 18.60027347771934
$ws.Range("J16").Value = 9.07697680127179
$ws.Range("K16").Value = 12.71339281754821
$ws.Range("L16").Value = 11.95732321920748
$ws.Range("N16").Value = 21.45592958680066
$ws.Range("B17").Value = 17.15780001756944
$ws.Range("C17").Value = 8.20994550012257
$ws.Range("E17").Value = 14.40567043590927
$ws.Range("F17").Value = 43.46723243898893
$ws.Range("G17").Value = 46.42773426595279
$ws.Range("H17").Value = 18.61467061487683
$ws.Range("J17").Value = 9.080595413016002
$ws.Range("K17").Value = 12.65230877044553
$ws.Range("L17").Value = 11.94754191814107
$ws.Range("N17").Value = 21.48201298779082
$ws.Range("B18").Value = 17.10633109930768
$ws.Range("C18").Value = 8.20342961659224
$ws.Range("E18").Value = 14.40314192322061
$ws.Range("F18").Value = 43.45257756691859
$ws.Range("G18").Value = 46.42273580778163
$ws.Range("H18").Value = 18.62323695433918
$ws.Range("J18").Value = 9.082702750712043
$ws.Range("K18").Value = 12.61728789627072
$ws.Range("L18").Value = 11.94208204183072
$ws.Range("N18").Value = 21.49719845725474
$ws.Range("B19").Value = 17.08892952750705
$ws.Range("C19").Value = 8.201216850154083
$ws.Range("E19").Value = 14.40231905839531
$ws.Range("F19").Value = 43.44777139101421
$ws.Range("G19").Value = 46.42123912105732
$ws.Range("H19").Value = 18.62618638596992
$ws.Range("J19").Value = 9.083420732229676
$ws.Range("K19").Value = 12.60545092128353
$ws.Range("L19").Value = 11.94026204948383
$ws.Range("N19").Value = 21.50237145827394
$ws.Range("B20").Value = 17.16733720450987
$ws.Range("C20").Value = 8.211148303640323
$ws.Range("E20").Value = 14.40615418260468
$ws.Range("F20").Value = 43.47001862868672
$ws.Range("G20").Value = 46.42875232343934
$ws.Range("H20").Value = 18.61310846432399
$ws.Range("J20").Value = 9.080207515385636
$ws.Range("K20").Value = 12.65879982616447
$ws.Range("L20").Value = 11.94856599200144
$ws.Range("N20").Value = 21.47921743473946
$ws.Range("B21").Value = 17.43227398201117
$ws.Range("C21").Value = 8.244054387353039
$ws.Range("E21").Value = 14.42134247721192
$ws.Range("F21").Value = 43.55555965376859
$ws.Range("G21").Value = 46.4676050355093
$ws.Range("H21").Value = 18.57258055976061
$ws.Range("J21").Value = 9.069721820837616
$ws.Range("K21").Value = 12.83931621960121
$ws.Range("L21").Value = 11.97844225516726
$ws.Range("N21").Value = 21.40360690197343
$ws.Range("B22").Value = 17.60626604139435
$ws.Range("C22").Value = 8.265212251908848
$ws.Range("E22").Value = 14.43297297824697
$ws.Range("F22").Value = 43.61946030915885
$ws.Range("G22").Value = 46.50313030102165
$ws.Range("H22").Value = 18.54863892381332
$ws.Range("J22").Value = 9.063099007687445
$ws.Range("K22").Value = 12.95806098612404
$ws.Range("L22").Value = 11.99942225132351
$ws.Range("N22").Value = 21.35581303893983
$ws.Range("B23").Value = 17.51334388756051
$ws.Range("C23").Value = 8.253951993044186
$ws.Range("E23").Value = 14.42661045220261
$ws.Range("F23").Value = 43.58462723977545
$ws.Range("G23").Value = 46.48324358275296
$ws.Range("H23").Value = 18.56118393539445
$ws.Range("J23").Value = 9.066612732603973
$ws.Range("K23").Value = 12.89462606176451
$ws.Range("L23").Value = 11.98809301753055
$ws.Range("N23").Value = 21.38117347048636
$ws.Range("B24").Value = 17.16302508879867
$ws.Range("C24").Value = 8.210604645357314
$ws.Range("E24").Value = 14.40593488267869
$ws.Range("F24").Value = 43.46875619420293
$ws.Range("G24").Value = 46.42828851914946
$ws.Range("H24").Value = 18.61381381164541
$ws.Range("J24").Value = 9.080382799934592
$ws.Range("K24").Value = 12.65586491453749
$ws.Range("L24").Value = 11.94810249858666
$ws.Range("N24").Value = 21.48048071301505
$ws.Range("B25").Value = 16.79054966949834
$ws.Range("C25").Value = 8.162376192228667
$ws.Range("E25").Value = 14.39102605503928
$ws.Range("F25").Value = 43.37839825709851
$ws.Range("G25").Value = 46.41256112257678
$ws.Range("H25").Value = 18.68149593512384
$ws.Range("J25").Value = 9.096231264419281
$ws.Range("K25").Value = 12.40278655380349
$ws.Range("L25").Value = 11.91132761316164
$ws.Range("N25").Value = 21.59460047149297

Write-Host "Applied 240 cell updates"
